$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.438.39"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "1.573.02"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.36"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3736"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.85"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3395"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.135"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.990"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.927"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "1.577.44"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001121"
$ws.Range("E17").Value = "  -1.14%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.05"
$ws.Range("E18").Value = "  +1.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06730"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.267"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("D24").Value = "22.449.36"
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.344"
$ws.Range("E25").Value = "  -3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.609"
$ws.Range("E26").Value = "  -4.35%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.44"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.013"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.65"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "1.753.66"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.054"
$ws.Range("E32").Value = "  +5.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.129"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.980"
$ws.Range("E34").Value = "  -1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.788"
$ws.Range("E35").Value = "  -2.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08393"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.380"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02456"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2287"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06515"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.26"
$ws.Range("E42").Value = "  -2.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6224"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.97"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5802"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.68"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.077"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.217"
$ws.Range("E50").Value = "  -7.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07321"
$ws.Range("E51").Value = "  -0.14%  "
